# backwardElimination.xlsx - refresh the statsmodels OLS "Summary()" text
# blocks (one per backward-elimination step / worksheet) to the timestamps
# from the re-run that produced the file actually committed:
#   Date: Sat, 28 Dec 2019 -> Sun, 29 Dec 2019   (every sheet)
#   Time: 20:59:35/20:59:36 -> 16:11:01/16:11:02 (per-sheet, see map below)
#
# Each worksheet is named after the number of remaining predictors in that
# backward-elimination step ("46" down to "18") and carries the full
# OLS Regression Results printout in cell B2.

$wb = $excel.ActiveWorkbook

# Per-sheet (tab name -> old/new "Time:" stamp) lookup, taken from the
# actual diff of the committed workbook.
$timeMap = @{
    "46" = @{ Old = "20:59:35"; New = "16:11:01" }
    "45" = @{ Old = "20:59:36"; New = "16:11:01" }
    "44" = @{ Old = "20:59:36"; New = "16:11:01" }
    "43" = @{ Old = "20:59:36"; New = "16:11:01" }
    "42" = @{ Old = "20:59:36"; New = "16:11:01" }
    "41" = @{ Old = "20:59:36"; New = "16:11:01" }
    "40" = @{ Old = "20:59:36"; New = "16:11:01" }
    "39" = @{ Old = "20:59:36"; New = "16:11:01" }
    "38" = @{ Old = "20:59:36"; New = "16:11:01" }
    "37" = @{ Old = "20:59:36"; New = "16:11:01" }
    "36" = @{ Old = "20:59:36"; New = "16:11:01" }
    "35" = @{ Old = "20:59:36"; New = "16:11:01" }
    "34" = @{ Old = "20:59:36"; New = "16:11:01" }
    "33" = @{ Old = "20:59:36"; New = "16:11:02" }
    "32" = @{ Old = "20:59:36"; New = "16:11:02" }
    "31" = @{ Old = "20:59:36"; New = "16:11:02" }
    "30" = @{ Old = "20:59:36"; New = "16:11:02" }
    "29" = @{ Old = "20:59:36"; New = "16:11:02" }
    "28" = @{ Old = "20:59:36"; New = "16:11:02" }
    "27" = @{ Old = "20:59:36"; New = "16:11:02" }
    "26" = @{ Old = "20:59:36"; New = "16:11:02" }
    "25" = @{ Old = "20:59:36"; New = "16:11:02" }
    "24" = @{ Old = "20:59:36"; New = "16:11:02" }
    "23" = @{ Old = "20:59:36"; New = "16:11:02" }
    "22" = @{ Old = "20:59:36"; New = "16:11:02" }
    "21" = @{ Old = "20:59:36"; New = "16:11:02" }
    "20" = @{ Old = "20:59:36"; New = "16:11:02" }
    "19" = @{ Old = "20:59:36"; New = "16:11:02" }
    "18" = @{ Old = "20:59:36"; New = "16:11:02" }
}

$oldDate = "Date:                Sat, 28 Dec 2019"
$newDate = "Date:                Sun, 29 Dec 2019"

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if (-not $timeMap.ContainsKey($name)) { continue }

    $cell = $ws.Range("B2")
    $text = $cell.Value2
    if (-not $text) { continue }
    if ($text.IndexOf("Dep. Variable") -lt 0) { continue }

    $oldTime = "Time:                        " + $timeMap[$name].Old
    $newTime = "Time:                        " + $timeMap[$name].New

    $updated = $text.Replace($oldDate, $newDate).Replace($oldTime, $newTime)

    if ($updated -ne $text) {
        $cell.Value2 = $updated
    }
}
